$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.045828333333333
$ws.Cells.Item(2, 8).Value = 3.137485
$ws.Cells.Item(2, 9).Value = 0.0274735881233672
$ws.Cells.Item(2, 10).Value = 0.0274735881233672
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 6.126464333333334
$ws.Cells.Item(2, 14).Value = 18.379393
$ws.Cells.Item(2, 15).Value = 0.1081098818071741
$ws.Cells.Item(2, 16).Value = 0.1081098818071741
$ws.Cells.Item(2, 17).Value = 6.407229982956112
$ws.Cells.Item(2, 18).Value = 57.665069846605
$ws.Cells.Item(2, 19).Value = 0.002970166364836211
$ws.Cells.Item(2, 20).Value = 0.002970166364836211

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.045828333333333
$ws.Cells.Item(3, 8).Value = 3.137485
$ws.Cells.Item(3, 9).Value = 0.0274735881233672
$ws.Cells.Item(3, 10).Value = 0.0274735881233672
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 35.73736866666667
$ws.Cells.Item(3, 14).Value = 107.212106
$ws.Cells.Item(3, 15).Value = 0.6306349784216607
$ws.Cells.Item(3, 16).Value = 0.6306349784216608
$ws.Cells.Item(3, 17).Value = 37.37515271037889
$ws.Cells.Item(3, 18).Value = 336.37637439341
$ws.Cells.Item(3, 19).Value = 0.01732580565334527
$ws.Cells.Item(3, 20).Value = 0.01732580565334527

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.045828333333333
$ws.Cells.Item(4, 8).Value = 3.137485
$ws.Cells.Item(4, 9).Value = 0.0274735881233672
$ws.Cells.Item(4, 10).Value = 0.0274735881233672
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 14.80503233333333
$ws.Cells.Item(4, 14).Value = 44.415097
$ws.Cells.Item(4, 15).Value = 0.2612551397711651
$ws.Cells.Item(4, 16).Value = 0.2612551397711651
$ws.Cells.Item(4, 17).Value = 15.48352229011611
$ws.Cells.Item(4, 18).Value = 139.351700611045
$ws.Cells.Item(4, 19).Value = 0.00717761610518572
$ws.Cells.Item(4, 20).Value = 0.00717761610518572

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 29.10096633333333
$ws.Cells.Item(5, 8).Value = 87.302899
$ws.Cells.Item(5, 9).Value = 0.7644734203038186
$ws.Cells.Item(5, 10).Value = 0.7644734203038187
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 6.126464333333334
$ws.Cells.Item(5, 14).Value = 18.379393
$ws.Cells.Item(5, 15).Value = 0.1081098818071741
$ws.Cells.Item(5, 16).Value = 0.1081098818071741
$ws.Cells.Item(5, 17).Value = 178.2860323067008
$ws.Cells.Item(5, 18).Value = 1604.574290760307
$ws.Cells.Item(5, 19).Value = 0.08264713111377198
$ws.Cells.Item(5, 20).Value = 0.08264713111377199

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 29.10096633333333
$ws.Cells.Item(6, 8).Value = 87.302899
$ws.Cells.Item(6, 9).Value = 0.7644734203038186
$ws.Cells.Item(6, 10).Value = 0.7644734203038187
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 35.73736866666667
$ws.Cells.Item(6, 14).Value = 107.212106
$ws.Cells.Item(6, 15).Value = 0.6306349784216607
$ws.Cells.Item(6, 16).Value = 0.6306349784216608
$ws.Cells.Item(6, 17).Value = 1039.991962410588
$ws.Cells.Item(6, 18).Value = 9359.927661695294
$ws.Cells.Item(6, 19).Value = 0.4821036789172318
$ws.Cells.Item(6, 20).Value = 0.4821036789172319

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 29.10096633333333
$ws.Cells.Item(7, 8).Value = 87.302899
$ws.Cells.Item(7, 9).Value = 0.7644734203038186
$ws.Cells.Item(7, 10).Value = 0.7644734203038187
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 14.80503233333333
$ws.Cells.Item(7, 14).Value = 44.415097
$ws.Cells.Item(7, 15).Value = 0.2612551397711651
$ws.Cells.Item(7, 16).Value = 0.2612551397711651
$ws.Cells.Item(7, 17).Value = 430.8407474962448
$ws.Cells.Item(7, 18).Value = 3877.566727466203
$ws.Cells.Item(7, 19).Value = 0.1997226102728148
$ws.Cells.Item(7, 20).Value = 0.1997226102728148

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 7.919887
$ws.Cells.Item(8, 8).Value = 23.759661
$ws.Cells.Item(8, 9).Value = 0.2080529915728142
$ws.Cells.Item(8, 10).Value = 0.2080529915728142
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 6.126464333333334
$ws.Cells.Item(8, 14).Value = 18.379393
$ws.Cells.Item(8, 15).Value = 0.1081098818071741
$ws.Cells.Item(8, 16).Value = 0.1081098818071741
$ws.Cells.Item(8, 17).Value = 48.52090522953034
$ws.Cells.Item(8, 18).Value = 436.688147065773
$ws.Cells.Item(8, 19).Value = 0.02249258432856593
$ws.Cells.Item(8, 20).Value = 0.02249258432856594

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 7.919887
$ws.Cells.Item(9, 8).Value = 23.759661
$ws.Cells.Item(9, 9).Value = 0.2080529915728142
$ws.Cells.Item(9, 10).Value = 0.2080529915728142
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 35.73736866666667
$ws.Cells.Item(9, 14).Value = 107.212106
$ws.Cells.Item(9, 15).Value = 0.6306349784216607
$ws.Cells.Item(9, 16).Value = 0.6306349784216608
$ws.Cells.Item(9, 17).Value = 283.0359215173407
$ws.Cells.Item(9, 18).Value = 2547.323293656066
$ws.Cells.Item(9, 19).Value = 0.1312054938510836
$ws.Cells.Item(9, 20).Value = 0.1312054938510837

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 7.919887
$ws.Cells.Item(10, 8).Value = 23.759661
$ws.Cells.Item(10, 9).Value = 0.2080529915728142
$ws.Cells.Item(10, 10).Value = 0.2080529915728142
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 14.80503233333333
$ws.Cells.Item(10, 14).Value = 44.415097
$ws.Cells.Item(10, 15).Value = 0.2612551397711651
$ws.Cells.Item(10, 16).Value = 0.2612551397711651
$ws.Cells.Item(10, 17).Value = 117.2541831113464
$ws.Cells.Item(10, 18).Value = 1055.287648002117
$ws.Cells.Item(10, 19).Value = 0.0543549133931646
$ws.Cells.Item(10, 20).Value = 0.05435491339316461
